$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": remove the row for client "FABIMP BENIGNO BRAVO S.A.S." (row 18) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(18).Delete()

# After the delete, the totals row (was 55, now 54) still shows the old denominator ("de 53").
# Update it to reflect the new number of clients ("de 52"), keeping each numerator unchanged.
$row1 = 54
$ws1.Cells.Item($row1, 3).Value2  = "1 de 52"
$ws1.Cells.Item($row1, 4).Value2  = "5 de 52"
$ws1.Cells.Item($row1, 5).Value2  = "3 de 52"
$ws1.Cells.Item($row1, 6).Value2  = "0 de 52"
$ws1.Cells.Item($row1, 7).Value2  = "0 de 52"
$ws1.Cells.Item($row1, 8).Value2  = "3 de 52"
$ws1.Cells.Item($row1, 9).Value2  = "9 de 52"
$ws1.Cells.Item($row1, 10).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 11).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 12).Value2 = "5 de 52"
$ws1.Cells.Item($row1, 13).Value2 = "11 de 52"
$ws1.Cells.Item($row1, 14).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 15).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 16).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 17).Value2 = "0 de 52"
$ws1.Cells.Item($row1, 18).Value2 = "0 de 52"

# --- Sheet "VENTA MENSUAL": remove the same client's row (row 18) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(18).Delete()

# The grand-total row (was 59, now 58) holds hard-coded sums, so subtract the removed row's
# contribution manually to keep the totals correct.
$row2 = 58
$ws2.Cells.Item($row2, 3).Value2 = 88156.04000000001
$ws2.Cells.Item($row2, 4).Value2 = 94831.00999999999
$ws2.Cells.Item($row2, 5).Value2 = 120013.04
$ws2.Cells.Item($row2, 6).Value2 = 28820.73
$ws2.Cells.Item($row2, 7).Value2 = 115950
